# Refresh cached market-price / profit figures on the Leve-profit sheets.
# (Mirrors a scheduled market-data sync; only H..N (price/profit) columns change.)
$wb = $excel.ActiveWorkbook

# ALC row 2: "Mercury Rising" (Quicksilver)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 204.25
$ws.Range("I2").Value = 204.25
$ws.Range("K2").Value = 204.25
$ws.Range("M2").Value = -91.25

# ALC row 41: "The Write Stuff" (Enchanted Mythril Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 255.90475
$ws.Range("I41").Value = 266.5
$ws.Range("J41").Value = 246.27272
$ws.Range("K41").Value = 266.5
$ws.Range("L41").Value = 246.27272
$ws.Range("M41").Value = 173.5
$ws.Range("N41").Value = -1126.27272

# ALC row 53: "No Accounting for Waste" (Enchanted Electrum Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 96.5
$ws.Range("I53").Value = 100.333336
$ws.Range("J53").Value = 94.85714
$ws.Range("K53").Value = 100.333336
$ws.Range("L53").Value = 94.85714
$ws.Range("M53").Value = 536.666664
$ws.Range("N53").Value = -1368.85714

# ALC row 54: "Arcane Arts for Dummies" (Book of Mythril)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 48999.75
$ws.Range("I54").Value = 54666.332
$ws.Range("K54").Value = 54666.332
$ws.Range("M54").Value = -54180.332

# ALC row 55: "A Real Smooth Move" (Lanolin)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 190
$ws.Range("I55").Value = 123.333336
$ws.Range("J55").Value = 240
$ws.Range("K55").Value = 123.333336
$ws.Range("L55").Value = 240
$ws.Range("M55").Value = 90.666664
$ws.Range("N55").Value = -668

# ALC row 64: "Forged from the Void" (Void Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3503.3
$ws.Range("I64").Value = 3257
$ws.Range("J64").Value = 3578.261
$ws.Range("K64").Value = 3257
$ws.Range("L64").Value = 3578.261
$ws.Range("M64").Value = -3009
$ws.Range("N64").Value = -4074.261

# ALC row 67: "Dodging the Draft (L)" (Void Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3503.3
$ws.Range("I67").Value = 3257
$ws.Range("J67").Value = 3578.261
$ws.Range("K67").Value = 3257
$ws.Range("L67").Value = 3578.261
$ws.Range("M67").Value = -2399
$ws.Range("N67").Value = -5294.261

# ALC row 74: "Adhesive of Antipathy" (Wing Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4119.091
$ws.Range("I74").Value = 3968.75
$ws.Range("J74").Value = 4205
$ws.Range("K74").Value = 3968.75
$ws.Range("L74").Value = 4205
$ws.Range("M74").Value = -3032.75
$ws.Range("N74").Value = -6077

# ALC row 77: "It's Gonna Grow Back (L)" (Wing Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4119.091
$ws.Range("I77").Value = 3968.75
$ws.Range("J77").Value = 4205
$ws.Range("K77").Value = 19843.75
$ws.Range("L77").Value = 21025
$ws.Range("M77").Value = -15163.75
$ws.Range("N77").Value = -30385

# ALC row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2566770
$ws.Range("I138").Value = 5715299
$ws.Range("J138").Value = 4013.7441
$ws.Range("K138").Value = 17145897
$ws.Range("L138").Value = 12041.2323
$ws.Range("M138").Value = -17140757
$ws.Range("N138").Value = -22321.2323

# ARM row 32: "Ingot We Trust" (Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11136.242
$ws.Range("I32").Value = 10887.018
$ws.Range("J32").Value = 14750
$ws.Range("K32").Value = 10887.018
$ws.Range("L32").Value = 14750
$ws.Range("M32").Value = -10600.018
$ws.Range("N32").Value = -15324

# BSM row 86: "Through Thick and Thin" (Adamantite Nugget)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 251491.5
$ws.Range("I86").Value = 1989.5
$ws.Range("J86").Value = 500993.5
$ws.Range("K86").Value = 1989.5
$ws.Range("L86").Value = 500993.5
$ws.Range("M86").Value = -866.5
$ws.Range("N86").Value = -503239.5

# BSM row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Adamantite Nugget)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 251491.5
$ws.Range("I89").Value = 1989.5
$ws.Range("J89").Value = 500993.5
$ws.Range("K89").Value = 9947.5
$ws.Range("L89").Value = 2504967.5
$ws.Range("M89").Value = -4331.5
$ws.Range("N89").Value = -2516199.5

# CRP row 62: "Splinter in the Sewers" (Cedar Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 127427.5
$ws.Range("I62").Value = 168903.33
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 168903.33
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -168279.33
$ws.Range("N62").Value = -4248

# CRP row 64: "Almost as Fun as Slingshotting Birds" (Cedar Longbow)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# CRP row 65: "The Lumber of Their Discontent (L)" (Cedar Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 127427.5
$ws.Range("I65").Value = 168903.33
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 844516.6499999999
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -841396.6499999999
$ws.Range("N65").Value = -21240

# CRP row 67: "Living Bow to Mouth (L)" (Cedar Longbow)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 30000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# CUL row 54: "Good Eats in Ishgard" (Salt Cod Puffs)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

# CUL row 55: "Pagan Pastries" (Pastry Fish)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 10277.4
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 12721.75
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 38165.25
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -38519.25

# GSM row 80: "Needs More Prayerbell" (Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3166.6667
$ws.Range("J80").Value = 3333.3333
$ws.Range("L80").Value = 3333.3333
$ws.Range("N80").Value = -5329.3333

# GSM row 83: "With a Noise That Reaches Heaven (L)" (Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3166.6667
$ws.Range("J83").Value = 3333.3333
$ws.Range("L83").Value = 16666.6665
$ws.Range("N83").Value = -26650.6665

# LTW row 25: "A Rush on Ringbands" (Hard Leather Ringbands)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 300335.66
$ws.Range("I25").Value = 300335.66
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 300335.66
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -300105.66
$ws.Range("N25").ClearContents()

# LTW row 59: "Fuss in Boots" (Raptorskin Thighboots)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41308

# LTW row 68: "You Could Say It's a Moving Target" (Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1988.8889
$ws.Range("I68").Value = 1580
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1580
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -831
$ws.Range("N68").Value = -3998

# LTW row 71: "They Call It Bloody Mary (L)" (Wyvern Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1988.8889
$ws.Range("I71").Value = 1580
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 7900
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -4156
$ws.Range("N71").Value = -19988

# LTW row 132: "Tenets of Tanning" (Silver Lobo Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5676.6665
$ws.Range("I132").Value = 5400.6665
$ws.Range("K132").Value = 16201.9995
$ws.Range("M132").Value = -13671.9995

# WVR row 114: "Hunting Season" (Pixie Cotton Hat of Striking)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# WVR row 132: "Comfy Cabins" (Snow Cotton Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1580.9565
$ws.Range("I132").Value = 1210
$ws.Range("J132").Value = 2428.8572
$ws.Range("K132").Value = 3630
$ws.Range("L132").Value = 7286.571599999999
$ws.Range("M132").Value = -1100
$ws.Range("N132").Value = -12346.5716
